# Commit: "metallurgy table - 1"
# Updates the RECIPES sheet status column for the Metallurgy Table recipes
# (re11/Pewter Bar, re12/Tin Ore, re13/Copper Ore) to "temporal", and adds the
# Tiny Investiture Spark material requirement to the Tin/Copper Ore recipes on
# the MATERIALS IN RECIPE sheet (also correcting the existing bar-vs-ore
# material references for those two recipes).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# RECIPES sheet: rows 12-14 (re11, re12, re13) -> status "temporal"
# ---------------------------------------------------------------------------
$recipes = $wb.Worksheets.Item("RECIPES")
$recipes.Range("G12").Value = "temporal"
$recipes.Range("G13").Value = "temporal"
$recipes.Range("G14").Value = "temporal"

# ---------------------------------------------------------------------------
# MATERIALS IN RECIPE sheet
# ---------------------------------------------------------------------------
$materials = $wb.Worksheets.Item("MATERIALS IN RECIPE")

# re12 (row 19) used Copper Ore - it should reference Copper Bar instead
$materials.Range("C19").Value = "Copper Bar"

# Insert a new row after row 19 for the re12 Tiny Investiture Spark requirement,
# pushing the old row 20 (re13 / Tin Ore) down to row 21.
$materials.Rows.Item(20).Insert()

$materials.Range("A20").Value = "re12"
$materials.Range("B20").Value = 2
$materials.Range("C20").Value = "Tiny Investiture Spark"
$materials.Range("D20").Value = 10

# Row 21 (previously row 20) referenced Tin Ore - it should reference Tin Bar instead
$materials.Range("C21").Value = "Tin Bar"

# Append the matching Tiny Investiture Spark requirement for re13
$materials.Range("A22").Value = "re13"
$materials.Range("B22").Value = 2
$materials.Range("C22").Value = "Tiny Investiture Spark"
$materials.Range("D22").Value = 10

# ---------------------------------------------------------------------------
# View selections (RECIPES must stay the active/selected sheet/tab, as in the
# original workbook - so set the MATERIALS IN RECIPE selection first, then
# finish by re-activating RECIPES).
# ---------------------------------------------------------------------------
$materials.Activate()
$materials.Range("A19:D20").Select()

$recipes.Activate()
$recipes.Range("G15").Select()
